$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dadosDeAcesso")

$ws.Range("A5").Value = "ID_0007"
$ws.Range("B5").Value = "André Automatizador"
$ws.Range("C5").Value = "sem email"
$ws.Range("D5").Value = "automacaoteste"

$ws.Range("C5").Select()
